# Add a "featureType" field to the single-line file format example.
#
# The worksheet holds three stacked mini-tables (rows 1-5, rows 6-8 and
# rows 9-11), each with its own little header row. Only the first table
# (the "feature" table, rows 1-5) actually gains the new column, inserted
# between "featureName" and "start". To get Excel to recompute the column
# width metadata (the <cols> bestFit widths) correctly we insert a real
# worksheet column (which shifts every row, including the other two
# tables) and then shift the unaffected tables' data back to where it
# was before the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "start" column (E). This shifts
# columns E..H to F..I for every row and correctly carries the bestFit
# column-width metadata along with it.
$ws.Columns.Item(5).Insert()

# --- Populate the new column for the first mini-table (rows 1-5) ---
$ws.Range("E1").Value = "featureType"
$ws.Range("E2").Value = "gene"
$ws.Range("E3").Value = "gene"
$ws.Range("E4").Value = "SNP"
$ws.Range("E5").Value = "SNP"

# --- The second mini-table (rows 6-8, the "homolog" table) should not
#     have shifted at all, so move its data back one column to the left ---
$ws.Range("F6").Cut($ws.Range("E6"))
$ws.Range("G6").Cut($ws.Range("F6"))
$ws.Range("H6").Cut($ws.Range("G6"))
$ws.Range("I6").Cut($ws.Range("H6"))
$ws.Range("I6").Clear()

$ws.Range("F7").Cut($ws.Range("E7"))
$ws.Range("G7").Cut($ws.Range("F7"))
$ws.Range("H7").Cut($ws.Range("G7"))
$ws.Range("H7").Clear()

$ws.Range("F8").Cut($ws.Range("E8"))
$ws.Range("G8").Cut($ws.Range("F8"))
$ws.Range("H8").Cut($ws.Range("G8"))
$ws.Range("H8").Clear()

# The third mini-table (rows 9-11) only has a trailing "comment" cell in
# its header row (row 9); the column insert already moved it from H9 to
# I9 correctly, so nothing else needs to be done there.

# --- Column widths: new column E gets a fixed custom width, and the
#     trailing "comment" column (now column I) gets a new fixed width ---
$ws.Columns.Item(5).ColumnWidth = 23.5
$ws.Columns.Item(9).ColumnWidth = 52.166666666666664

# --- Selection moved as part of the edit ---
$ws.Range("F19").Select()
